$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2: M2 -> Cd28/Cd86 -> M2 (updated values)
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Cd28"
$ws.Range("C2").Value = "Cd86"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.543204
$ws.Range("H2").Value = 22.629612
$ws.Range("I2").Value = 0.9718881576768906
$ws.Range("J2").Value = 0.9718881576768906
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 84.02338933333334
$ws.Range("N2").Value = 252.070168
$ws.Range("O2").Value = 0.9856776291044543
$ws.Range("P2").Value = 0.9856776291044543
$ws.Range("Q2").Value = 633.8055665127574
$ws.Range("R2").Value = 5704.250098614817
$ws.Range("S2").Value = 0.9579684150136536
$ws.Range("T2").Value = 0.9579684150136536

# New row 3: M2 -> Cd28/Cd86 -> sCs
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Cd28"
$ws.Range("C3").Value = "Cd86"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.543204
$ws.Range("H3").Value = 22.629612
$ws.Range("I3").Value = 0.9718881576768906
$ws.Range("J3").Value = 0.9718881576768906
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 1.220900333333333
$ws.Range("N3").Value = 3.662701
$ws.Range("O3").Value = 0.0143223708955457
$ws.Range("P3").Value = 0.01432237089554569
$ws.Range("Q3").Value = 9.209500278001334
$ws.Range("R3").Value = 82.88550250201202
$ws.Range("S3").Value = 0.01391974266323702
$ws.Range("T3").Value = 0.01391974266323702

# New row 4: sCs -> Cd28/Cd86 -> M2
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Cd28"
$ws.Range("C4").Value = "Cd86"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.218187
$ws.Range("H4").Value = 0.6545609999999999
$ws.Range("I4").Value = 0.02811184232310935
$ws.Range("J4").Value = 0.02811184232310934
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 84.02338933333334
$ws.Range("N4").Value = 252.070168
$ws.Range("O4").Value = 0.9856776291044543
$ws.Range("P4").Value = 0.9856776291044543
$ws.Range("Q4").Value = 18.332811248472
$ws.Range("R4").Value = 164.995301236248
$ws.Range("S4").Value = 0.02770921409080068
$ws.Range("T4").Value = 0.02770921409080067

# New row 5: sCs -> Cd28/Cd86 -> sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Cd28"
$ws.Range("C5").Value = "Cd86"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.218187
$ws.Range("H5").Value = 0.6545609999999999
$ws.Range("I5").Value = 0.02811184232310935
$ws.Range("J5").Value = 0.02811184232310934
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 1.220900333333333
$ws.Range("N5").Value = 3.662701
$ws.Range("O5").Value = 0.0143223708955457
$ws.Range("P5").Value = 0.01432237089554569
$ws.Range("Q5").Value = 0.266384581029
$ws.Range("R5").Value = 2.397461229261
$ws.Range("S5").Value = 0.000402628232308671
$ws.Range("T5").Value = 0.0004026282323086709
